$d = $word.ActiveDocument

# 1. Remove the stray "_GoBack" bookmark that currently sits at the end of
#    the "Actual auction resolve..." paragraph (it is being relocated).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Split the run "...folder. Handle function reads periodically all
#    auctions saved to system and check auction status and auction ending
#    time. If auction status " into three runs:
#      a) "...all auctions saved "
#      b) "to system and check auction status and auction ending time"
#      c) ". If auction status "
#    First split point: right after "all auctions saved ".
$rng = $d.Content
$rng.Find.Execute("all auctions saved ", $false)
$rng.Collapse(0)
$splitPos1 = $rng.End
$tmpRange = $d.Range($splitPos1, $splitPos1)
$d.Bookmarks.Add("TEMP_SPLIT_MARKER", $tmpRange)
$d.Bookmarks("TEMP_SPLIT_MARKER").Delete()

# Second split point: right after "...auction ending time" -- this is also
# where the relocated "_GoBack" bookmark belongs.
$rng2 = $d.Content
$rng2.Find.Execute("to system and check auction status and auction ending time", $false)
$rng2.Collapse(0)
$splitPos2 = $rng2.End
$goBackRange = $d.Range($splitPos2, $splitPos2)
$d.Bookmarks.Add("_GoBack", $goBackRange)
